# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Clean up the "ODI Batting Extra" sheet: a handful of cells were written
#    with an explicit-but-empty string; the refreshed export drops those
#    cells entirely instead of leaving an empty placeholder behind.
# 2) Add the new "ODI Bowling Extra" sheet (mirrors "ODI Batting Extra")
#    holding the MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL scrape for each match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "ODI Batting Extra" -- drop the empty placeholder cells
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ODI Batting Extra")

$ws4.Range("C3:E3").ClearContents()
$ws4.Range("B4:E4").ClearContents()
$ws4.Range("C5:E5").ClearContents()
$ws4.Range("B6:E6").ClearContents()
$ws4.Range("B7:E7").ClearContents()
$ws4.Range("B9:E9").ClearContents()
$ws4.Range("E12").ClearContents()
$ws4.Range("C14:E14").ClearContents()
$ws4.Range("B18:E18").ClearContents()
$ws4.Range("C21:E21").ClearContents()

# ---------------------------------------------------------------------
# 2) Add "ODI Bowling Extra" as the new last sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws5.Name = "ODI Bowling Extra"

# Header row, styled like the other "*_Extra" sheets (bold, centered, boxed)
$header = $ws5.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

$ws5.Range("A1").Value = "MATCH_CODE"
$ws5.Range("B1").Value = "MAIDEN_OVERS"
$ws5.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Data rows -- values are written as text (leading "'" forces text, matching
# the source sheet where MATCH_CODE / MAIDEN_OVERS / the percent strings are
# all stored as strings rather than numbers)
$ws5.Range("A2").Value = "'4251"
$ws5.Range("B2").Value = "'0"
$ws5.Range("A3").Value = "'4252"
$ws5.Range("B3").Value = "'0"
$ws5.Range("C3").Value = "'10.00%"
$ws5.Range("A4").Value = "'4286"
$ws5.Range("B4").Value = "'1"
$ws5.Range("C4").Value = "'20.00%"
$ws5.Range("A5").Value = "'4295"
$ws5.Range("A6").Value = "'4296"
$ws5.Range("B6").Value = "'0"
$ws5.Range("A7").Value = "'4307"
$ws5.Range("A8").Value = "'4311"
$ws5.Range("A9").Value = "'4314"
$ws5.Range("B9").Value = "'0"
$ws5.Range("C9").Value = "'20.00%"
$ws5.Range("A10").Value = "'4325"
$ws5.Range("A11").Value = "'4335"
$ws5.Range("B11").Value = "'0"
$ws5.Range("C11").Value = "'10.00%"
$ws5.Range("A12").Value = "'4345"
$ws5.Range("B12").Value = "'0"
$ws5.Range("A13").Value = "'4349"
$ws5.Range("B13").Value = "'0"
$ws5.Range("C13").Value = "'30.00%"
$ws5.Range("A14").Value = "'4416"
$ws5.Range("B14").Value = "'0"
$ws5.Range("C14").Value = "'30.00%"
$ws5.Range("A15").Value = "'4420"
$ws5.Range("B15").Value = "'0"
$ws5.Range("C15").Value = "'40.00%"
$ws5.Range("A16").Value = "'4447"
$ws5.Range("B16").Value = "'0"
$ws5.Range("C16").Value = "'30.00%"
$ws5.Range("A17").Value = "'4453"
$ws5.Range("B17").Value = "'0"
$ws5.Range("A18").Value = "'4463"
$ws5.Range("B18").Value = "'0"
$ws5.Range("C18").Value = "'20.00%"
$ws5.Range("A19").Value = "'4477"
$ws5.Range("B19").Value = "'0"
$ws5.Range("C19").Value = "'10.00%"
$ws5.Range("A20").Value = "'4479"
$ws5.Range("B20").Value = "'0"
$ws5.Range("C20").Value = "'10.00%"
$ws5.Range("A21").Value = "'4481"
$ws5.Range("B21").Value = "'0"
$ws5.Range("C21").Value = "'30.00%"

$ws5.Range("A1").Select()

Write-Host "Sheets now:" ($wb.Worksheets.Count)
